# Weekly update: two new price observations (date serial 44610) were
# recorded for Cilantro at Vega Central Mapocho de Santiago, inserted right
# after the existing pair dated 44312 (rows 488-489), pushing every
# subsequent row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 490 (shifts old rows 490:516 -> 492:518)
$ws.Rows.Item(490).Insert()
$ws.Rows.Item(490).Insert()

# New row 490 - $/caja 36 atados observation
$ws.Cells.Item(490, 1).Value = 9
$ws.Cells.Item(490, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(490, 3).Value = "Metropolitana"
$ws.Cells.Item(490, 4).Value = 44610
$ws.Cells.Item(490, 5).Value = 13
$ws.Cells.Item(490, 6).Value = 100112040
$ws.Cells.Item(490, 7).Value = "Cilantro"
$ws.Cells.Item(490, 8).Value = "Sin especificar"
$ws.Cells.Item(490, 9).Value = "Primera"
$ws.Cells.Item(490, 10).Value = 52
$ws.Cells.Item(490, 11).Value = 9000
$ws.Cells.Item(490, 12).Value = 9000
$ws.Cells.Item(490, 13).Value = 9000
$ws.Cells.Item(490, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(490, 15).Value = "Región Metropolitana"
$ws.Cells.Item(490, 16).Value = 250
$ws.Cells.Item(490, 17).Value = 36
$ws.Cells.Item(490, 18).Value = "Hortaliza"

# New row 491 - $/docena de atados observation
$ws.Cells.Item(491, 1).Value = 9
$ws.Cells.Item(491, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(491, 3).Value = "Metropolitana"
$ws.Cells.Item(491, 4).Value = 44610
$ws.Cells.Item(491, 5).Value = 13
$ws.Cells.Item(491, 6).Value = 100112040
$ws.Cells.Item(491, 7).Value = "Cilantro"
$ws.Cells.Item(491, 8).Value = "Sin especificar"
$ws.Cells.Item(491, 9).Value = "Primera"
$ws.Cells.Item(491, 10).Value = 160
$ws.Cells.Item(491, 11).Value = 16000
$ws.Cells.Item(491, 12).Value = 18000
$ws.Cells.Item(491, 13).Value = 17000
$ws.Cells.Item(491, 14).Value = "`$/docena de atados"
$ws.Cells.Item(491, 15).Value = "Región Metropolitana"
$ws.Cells.Item(491, 16).Value = 5667
$ws.Cells.Item(491, 17).Value = 3
$ws.Cells.Item(491, 18).Value = "Hortaliza"
